$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets "fresh" data: the CO number formerly on row 4 moves up to row 2.
$ws.Range("A2").Value = "3013696594"

# The stale rows 3-6 are wiped out completely (content + formatting), leaving
# an empty gap of rows 3-9.
$ws.Range("A3:O9").Clear()

# The remaining two "fresh" CO numbers land on rows 10 and 11, replicating the
# same row layout used by the other data rows.
$dataCols    = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")
$leftCols    = @("A","B","E","F","H","K","L","M","N")   # style index 1 (left/top)
# everything else (C,D,G,I,J,O) uses style index 3 (right/top)

$row10 = @{
    "A" = "3013696595"; "B" = "TA5ACMFTWLC"; "C" = "101017"; "D" = "1.00";
    "E" = "EA";          "F" = "20";          "G" = "101017"; "H" = "";
    "I" = "";             "J" = "";            "K" = "A2";    "L" = " 7";
    "M" = "001";          "N" = "";            "O" = ""
}
$row11 = @{
    "A" = "3013696596"; "B" = "TA5ACMFTWLC"; "C" = "101017"; "D" = "1.00";
    "E" = "EA";          "F" = "20";          "G" = "101017"; "H" = "";
    "I" = "";             "J" = "";            "K" = "A2";    "L" = " 7";
    "M" = "001";          "N" = "";            "O" = ""
}

foreach ($r in @(10, 11)) {
    $rowData = if ($r -eq 10) { $row10 } else { $row11 }
    foreach ($c in $dataCols) {
        $cell = $ws.Range($c + $r)
        $cell.Value = $rowData[$c]
        $cell.NumberFormat = "@"
        $cell.VerticalAlignment = -4160   # xlTop
        if ($leftCols -contains $c) {
            $cell.HorizontalAlignment = -4131   # xlLeft
        } else {
            $cell.HorizontalAlignment = -4152   # xlRight
        }
    }
}

# Update the view's selection to match the post-edit state.
$ws.Range("A6:XFD6").Select()
